$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark every remaining "n" (not-yet-updated) part status as "y" (finished
# updating parts in schematic).
$statusRows = @(2,4,5,6,7,8,9,10,11,12,13,14,16,17,21,22,23,24,25,26,27,28,29,34)
foreach ($r in $statusRows) {
    $ws.Range("F$r").Value = "y"
}

# The note about the 150060GS75000 pricing is no longer relevant - clear it
# (keeping the existing cell formatting).
$ws.Range("G16").ClearContents()

# Update the last active selection left by the author.
$ws.Range("G15").Select()

Write-Output "done"
